$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target range grows from A1:I4 to A1:K6 (two new columns inserted for ownTeam/oppTeam,
# existing rows reordered/augmented, and two new match rows appended).
# Force text storage (matches the source file, which stores every value as text,
# including numeric-looking ones) without leaving a stray number-format style behind.
$fullRange = $ws.Range("A1:K6")
$fullRange.NumberFormat = "@"

$data = @(
,@("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")
,@(" Abu Dhabi"," October 30 2020","Royals won by 7 wickets (with 15 balls remaining)","Kings XI Punjab","Rajasthan Royals","Deepak Hooda ","1","1","0","0","100.00")
,@(" Dubai (DSC)"," October 24 2020","Kings XI won by 12 runs","Kings XI Punjab","Sunrisers Hyderabad","Deepak Hooda ","0","2","0","0","0.00")
,@(" Abu Dhabi"," November 01 2020","Super Kings won by 9 wickets (with 7 balls remaining)","Kings XI Punjab","Chennai Super Kings","Deepak Hooda ","62","30","3","4","206.66")
,@(" Dubai (DSC)"," October 20 2020","Kings XI won by 5 wickets (with 6 balls remaining)","Kings XI Punjab","Delhi Capitals","Deepak Hooda ","15","22","1","0","68.18")
,@(" Dubai (DSC)"," October 18 2020","Match tied (Kings XI won the one-over eliminator)","Kings XI Punjab","Mumbai Indians","Deepak Hooda ","23","16","1","1","143.75")
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# Restore the default "Normal" style on the whole range so no extra number format
# lingers on the cells themselves (the text-storage behaviour above is keyed off the
# assignment, not off the cell's final display format).
$fullRange.Style = "Normal"

